$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "Светлый, Калининградская область" "Светлый"
Replace-Text "генеральный директор Котлярчук О. Ю." "генеральный директор Котлярчук О. Е."
Replace-Text "Подтверждение Свидетельства об одобрении технологического процесса сварки (без испытаний) - 1 шт." "Первичная аттестация сваршиков - 2 чел."
Replace-Text "Свидетельство ф. 7.1.30 № 24.42.03.23414.121 от 29.04.2024" "Свидетельство ф. 7.1.30 № 24.42.03.00414.121 от 29.04.2024"
Replace-Text "О. Ю. Котлярчук" "О. Е. Котлярчук"
